$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 2 data: a new "quete" entry ("Les n nombres") ---
# (cells are populated in the same order the original author typed them in,
#  so that the shared-strings table comes out in the same order)

# Nom_Quete
$ws.Range("B2").Value = "Les n nombres"

# Rappel_de_la_quete
$ws.Range("L2").Value = "Ecrivez un programme qui crée une liste vide nommée L et qui ajoute ensuite à cette liste n nombres entiers, compris entre 1 et 10, saisis par l'utilisateur. `nFinalement, le programme affichera toute la liste ainsi que chaque élément de la liste séparément"

# input_exemple_1 / 2 / 3
$ws.Range("F2").Value = "5`n1`n2`n3`n4`n5"
$ws.Range("G2").Value = "2`n6`n9"
$ws.Range("H2").Value = "8`n4`n4`n4`n4`n6`n6`n4`n5"

# output_attendu_1 / 2 / 3
$ws.Range("I2").Value = "[1, 2, 3, 4, 5]`nL[0] = 1`nL[1] = 2`nL[2] = 3`nL[3] = 4`nL[4] = 5"
$ws.Range("J2").Value = "[6, 9]`nL[0] = 6`nL[1] = 9"
$ws.Range("K2").Value = "[4, 4, 4, 4, 6, 6, 4, 5]`nL[0] = 4`nL[1] = 4`nL[2] = 4`nL[3] = 4`nL[4] = 6`nL[5] = 6`nL[6] = 4`nL[7] = 5"

# Chap
$ws.Range("A2").Value = 4
# ID quete
$ws.Range("C2").Value = 0
# Difficulte
$ws.Range("D2").Value = 1

# n_test
$ws.Range("Q2").Value = 3

# input_test_1 / 2 / 3 (same content as input_exemple_1/2/3)
$ws.Range("R2").Value = "5`n1`n2`n3`n4`n5"
$ws.Range("T2").Value = "2`n6`n9"
$ws.Range("V2").Value = "8`n4`n4`n4`n4`n6`n6`n4`n5"

# output_test_1 / 2 / 3 (same content as output_attendu_1/2/3)
$ws.Range("S2").Value = "[1, 2, 3, 4, 5]`nL[0] = 1`nL[1] = 2`nL[2] = 3`nL[3] = 4`nL[4] = 5"
$ws.Range("U2").Value = "[6, 9]`nL[0] = 6`nL[1] = 9"
$ws.Range("W2").Value = "[4, 4, 4, 4, 6, 6, 4, 5]`nL[0] = 4`nL[1] = 4`nL[2] = 4`nL[3] = 4`nL[4] = 6`nL[5] = 6`nL[6] = 4`nL[7] = 5"

# --- Column B width ---
$ws.Columns.Item(2).ColumnWidth = 13.33

# --- View: zoom + selection on J3 ---
$excel.ActiveWindow.Zoom = 77
$ws.Range("J3").Select()
